# Auto-generated edit script: updates Jenova_Profits leve-profit computation columns
# (currentAveragePrice/currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# per the authoritative diff, across the 8 crafting-class worksheets.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets("ALC")
$ws.Range("H29").Value = 1390
$ws.Range("J29").Value = 3600
$ws.Range("L29").Value = 10800
$ws.Range("N29").Value = -11362
$ws.Range("H70").Value = 59964.47
$ws.Range("I70").Value = 815
$ws.Range("K70").Value = 2445
$ws.Range("M70").Value = -2175
$ws.Range("H73").Value = 59964.47
$ws.Range("I73").Value = 815
$ws.Range("K73").Value = 2445
$ws.Range("M73").Value = -1509
$ws.Range("H88").Value = 3946.9546
$ws.Range("J88").Value = 3887.7144
$ws.Range("L88").Value = 3887.7144
$ws.Range("N88").Value = -4699.7144
$ws.Range("H91").Value = 3946.9546
$ws.Range("J91").Value = 3887.7144
$ws.Range("L91").Value = 3887.7144
$ws.Range("N91").Value = -6695.7144
$ws.Range("H116").Value = 4214.2856
$ws.Range("I116").Value = 3900
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 3900
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -458
$ws.Range("N116").Value = -11884
$ws.Range("H133").Value = 54580.332
$ws.Range("J133").Value = 54580.332
$ws.Range("L133").Value = 54580.332
$ws.Range("N133").Value = -64700.332
$ws.Range("H137").Value = 3362.9583
$ws.Range("I137").Value = 2288.6667
$ws.Range("J137").Value = 4437.25
$ws.Range("K137").Value = 6866.000100000001
$ws.Range("L137").Value = 13311.75
$ws.Range("M137").Value = -4316.000100000001
$ws.Range("N137").Value = -18411.75

# ---- ARM ----
$ws = $wb.Worksheets("ARM")
$ws.Range("H74").Value = 16436.6
$ws.Range("I74").Value = 20223.334
$ws.Range("K74").Value = 20223.334
$ws.Range("M74").Value = -19349.334
$ws.Range("H77").Value = 16436.6
$ws.Range("I77").Value = 20223.334
$ws.Range("K77").Value = 101116.67
$ws.Range("M77").Value = -96748.67
$ws.Range("H102").Value = 1562.1
$ws.Range("I102").Value = 1291.6
$ws.Range("K102").Value = 1291.6
$ws.Range("M102").Value = 330.4000000000001
$ws.Range("H122").Value = 2911.3171
$ws.Range("I122").Value = 2018.931
$ws.Range("K122").Value = 6056.793
$ws.Range("M122").Value = -3606.793
$ws.Range("H128").Value = 73000
$ws.Range("J128").Value = 73000
$ws.Range("L128").Value = 73000
$ws.Range("N128").Value = -82960

# ---- BSM ----
$ws = $wb.Worksheets("BSM")
$ws.Range("H94").Value = 1082.88
$ws.Range("I94").Value = 796.619
$ws.Range("J94").Value = 2585.75
$ws.Range("K94").Value = 796.619
$ws.Range("L94").Value = 2585.75
$ws.Range("M94").Value = -345.619
$ws.Range("N94").Value = -3487.75

# ---- CRP ----
$ws = $wb.Worksheets("CRP")
$ws.Range("H31").Value = 3835.2144
$ws.Range("J31").Value = 6649
$ws.Range("L31").Value = 6649
$ws.Range("N31").Value = -7239
$ws.Range("H34").Value = 3835.2144
$ws.Range("J34").Value = 6649
$ws.Range("L34").Value = 6649
$ws.Range("N34").Value = -7053
$ws.Range("H122").Value = 2874.5454
$ws.Range("I122").Value = 1825.75
$ws.Range("J122").Value = 5671.3335
$ws.Range("K122").Value = 5477.25
$ws.Range("L122").Value = 17014.0005
$ws.Range("M122").Value = -3027.25
$ws.Range("N122").Value = -21914.0005
$ws.Range("H132").Value = 3745.611
$ws.Range("I132").Value = 3109
$ws.Range("K132").Value = 9327
$ws.Range("M132").Value = -6797

# ---- CUL ----
$ws = $wb.Worksheets("CUL")
$ws.Range("H2").Value = 12.045455
$ws.Range("I2").Value = 10.769231
$ws.Range("K2").Value = 64.615386
$ws.Range("M2").Value = 48.384614
$ws.Range("H5").Value = 90253.44500000001
$ws.Range("I5").Value = 100285.125
$ws.Range("K5").Value = 300855.375
$ws.Range("M5").Value = -300743.375
$ws.Range("H38").Value = 29.625
$ws.Range("I38").Value = 50
$ws.Range("J38").Value = 22.833334
$ws.Range("K38").Value = 150
$ws.Range("L38").Value = 68.50000199999999
$ws.Range("M38").Value = 197
$ws.Range("N38").Value = -762.500002
$ws.Range("H107").Value = 3169.1052
$ws.Range("J107").Value = 3615.5
$ws.Range("L107").Value = 10846.5
$ws.Range("N107").Value = -14686.5
$ws.Range("H110").Value = 2583.3333
$ws.Range("I110").Value = 2583.3333
$ws.Range("K110").Value = 7749.999899999999
$ws.Range("M110").Value = -3659.999899999999
$ws.Range("H122").Value = 84483.164
$ws.Range("J122").Value = 92072.55
$ws.Range("L122").Value = 828652.9500000001
$ws.Range("N122").Value = -833552.9500000001
$ws.Range("H135").Value = 90253.44500000001
$ws.Range("I135").Value = 100285.125
$ws.Range("K135").Value = 902566.125
$ws.Range("M135").Value = -900031.125
$ws.Range("H136").Value = 4811.857
$ws.Range("I136").Value = 4164
$ws.Range("J136").Value = 8699
$ws.Range("K136").Value = 12492
$ws.Range("L136").Value = 26097
$ws.Range("M136").Value = -7392
$ws.Range("N136").Value = -36297

# ---- GSM ----
$ws = $wb.Worksheets("GSM")
$ws.Range("H97").Value = 1314
$ws.Range("I97").Value = 1504.8572
$ws.Range("K97").Value = 1504.8572
$ws.Range("M97").Value = -1008.8572
$ws.Range("H102").Value = 5920
$ws.Range("I102").Value = 6873.273
$ws.Range("K102").Value = 6873.273
$ws.Range("M102").Value = -5251.273

# ---- LTW ----
$ws = $wb.Worksheets("LTW")
$ws.Range("H7").Value = 441541.97
$ws.Range("I7").Value = 840606.75
$ws.Range("K7").Value = 840606.75
$ws.Range("M7").Value = -840494.75
$ws.Range("H16").Value = 1408.3334
$ws.Range("I16").Value = 1290
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1290
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -1120
$ws.Range("N16").Value = -2340
$ws.Range("H20").Value = 50000
$ws.Range("J20").Value = 50000
$ws.Range("L20").Value = 50000
$ws.Range("N20").Value = -50452
$ws.Range("H122").Value = 754779.0600000001
$ws.Range("I122").Value = 670134.25
$ws.Range("J122").Value = 852446.1
$ws.Range("K122").Value = 2010402.75
$ws.Range("L122").Value = 2557338.3
$ws.Range("M122").Value = -2007952.75
$ws.Range("N122").Value = -2562238.3
$ws.Range("H126").Value = 441541.97
$ws.Range("I126").Value = 840606.75
$ws.Range("K126").Value = 2521820.25
$ws.Range("M126").Value = -2519350.25
$ws.Range("H132").Value = 4660.25
$ws.Range("I132").Value = 3014.7144
$ws.Range("J132").Value = 8499.833000000001
$ws.Range("K132").Value = 9044.143199999999
$ws.Range("L132").Value = 25499.499
$ws.Range("M132").Value = -6514.143199999999
$ws.Range("N132").Value = -30559.499
$ws.Range("H139").Value = 47905
$ws.Range("J139").Value = 47905
$ws.Range("L139").Value = 47905
$ws.Range("N139").Value = -58185

# ---- WVR ----
$ws = $wb.Worksheets("WVR")
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H81").Value = 18438.133
$ws.Range("I81").Value = 19437.857
$ws.Range("J81").Value = 4442
$ws.Range("K81").Value = 38875.714
$ws.Range("L81").Value = 8884
$ws.Range("M81").Value = -37814.714
$ws.Range("N81").Value = -11006
$ws.Range("H84").Value = 18438.133
$ws.Range("I84").Value = 19437.857
$ws.Range("J84").Value = 4442
$ws.Range("K84").Value = 194378.57
$ws.Range("L84").Value = 44420
$ws.Range("M84").Value = -189074.57
$ws.Range("N84").Value = -55028
$ws.Range("H132").Value = 36718.418
$ws.Range("I132").Value = 2758.6924
$ws.Range("J132").Value = 61244.89
$ws.Range("K132").Value = 8276.0772
$ws.Range("L132").Value = 183734.67
$ws.Range("M132").Value = -5746.0772
$ws.Range("N132").Value = -188794.67
